# Statistical Analysis Results - add MADRS score stats for depressed
# patients (Condition group), matching the "Included statistical analysis
# results for MADRS scores in depressed patients" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rows 3:5 (the "Number of days" table) pick up the same
#     border+center style used by the new MADRS rows below. ---
$ws.Range("B3:H5").HorizontalAlignment = -4108
$ws.Range("B3:H5").Borders.LineStyle = 1

# --- New "Start (madrs1)" row (row 10) ---
$ws.Range("B10").Value = 22.7
$ws.Range("C10").Value = "24, 26"
$ws.Range("D10").Value = 24
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 29
$ws.Range("G10").Value = 4.8
$ws.Range("H10").Value = 23

# --- New "End (madrs2)" row (row 11) ---
$ws.Range("B11").Value = 20
$ws.Range("C11").Value = 21
$ws.Range("D11").Value = 21
$ws.Range("E11").Value = 11
$ws.Range("F11").Value = 28
$ws.Range("G11").Value = 4.7
$ws.Range("H11").Value = 22.4

# Apply the same border + centered style to the newly-populated rows.
$ws.Range("B10:H11").HorizontalAlignment = -4108
$ws.Range("B10:H11").Borders.LineStyle = 1

# --- The old "Difference" row (row 12) no longer holds any data; clear
#     it out entirely (the label and the empty value cells), including
#     the border it used to carry. ---
$ws.Range("A12:H12").ClearContents()
$ws.Range("A12:H12").Borders.LineStyle = -4142
